$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.826.39'
$ws.Range("E2").Value = '  -0.89%  '
$ws.Range("D3").Value = '1.940.10'
$ws.Range("E3").Value = '  -0.61%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.89'
$ws.Range("E5").Value = '  -1.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2941'
$ws.Range("E8").Value = '  -0.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06919'
$ws.Range("E9").Value = '  +1.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.41'
$ws.Range("E10").Value = '  +2.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '105.68'
$ws.Range("E11").Value = '  -0.55%  '
$ws.Range("D12").Value = '1.941.80'
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07723'
$ws.Range("E13").Value = '  +0.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.359'
$ws.Range("E14").Value = '  -0.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6993'
$ws.Range("E15").Value = '  -1.63%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '273.65'
$ws.Range("E16").Value = '  -4.44%  '
$ws.Range("D17").Value = '30.824.99'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007739'
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.10'
$ws.Range("E19").Value = '  -0.62%  '
$ws.Range("D20").Value = '2.203.64'
$ws.Range("E20").Value = '  +0.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.530'
$ws.Range("E22").Value = '  +0.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.566'
$ws.Range("E24").Value = '  -0.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.729'
$ws.Range("E25").Value = '  -1.71%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.76'
$ws.Range("E26").Value = '  -0.94%  '
$ws.Range("E27").Value = '  -1.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.162'
$ws.Range("E28").Value = '  -1.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1039'
$ws.Range("E29").Value = '  -1.02%  '
$ws.Range("E30").Value = '  -3.31%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.577'
$ws.Range("E31").Value = '  -2.97%  '
$ws.Range("E32").Value = '  -2.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.368'
$ws.Range("E33").Value = '  -1.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04858'
$ws.Range("E34").Value = '  -2.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7566'
$ws.Range("E35").Value = '  -0.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.156'
$ws.Range("E36").Value = '  -0.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.000'
$ws.Range("E37").Value = '  +0.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.735'
$ws.Range("E38").Value = '  +0.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01995'
$ws.Range("E39").Value = '  -2.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.659'
$ws.Range("E40").Value = '  -1.60%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.489'
$ws.Range("E41").Value = '  +1.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '77.47'
$ws.Range("E42").Value = '  +6.94%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.091'
$ws.Range("E43").Value = '  -2.38%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9037'
$ws.Range("E44").Value = '  +2.80%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4412'
$ws.Range("E45").Value = '  -1.30%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '108.07'
$ws.Range("E46").Value = '  -1.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9987'
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.781'
$ws.Range("E48").Value = '  +4.43%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '988.35'
$ws.Range("E49").Value = '  +0.90%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1247'
$ws.Range("E50").Value = '  -2.40%  '
$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '36.12'
$ws.Range("E51").Value = '  +0.63%  '
